# Updates the cryptos list (rows 2-51) to the latest scraped values.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "40.019.13"; E = "  +1.44%  " },
    @{ Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "2.194.83"; E = "  +1.47%  " },
    @{ Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "1.00"; E = "  +0.26%  " },
    @{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "228.00"; E = "  -0.54%  " },
    @{ Row = 6; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "0.630"; E = "  +1.26%  " },
    @{ Row = 7; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "63.63"; E = "  +0.63%  " },
    @{ Row = 8; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "1.00"; E = "  +0.12%  " },
    @{ Row = 9; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.396"; E = "  +0.12%  " },
    @{ Row = 10; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.0859"; E = "  -0.56%  " },
    @{ Row = 11; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.104"; E = "  +0.19%  " },
    @{ Row = 12; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "2.528.38"; E = "  +1.79%  " },
    @{ Row = 13; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "15.85"; E = "  -0.98%  " },
    @{ Row = 14; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "22.04"; E = "  -0.72%  " },
    @{ Row = 15; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "0.816"; E = "  -0.05%  " },
    @{ Row = 16; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "5.56"; E = "  -0.40%  " },
    @{ Row = 17; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "2.211.65"; E = "  +2.11%  " },
    @{ Row = 18; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "39.994.30"; E = "  +1.41%  " },
    @{ Row = 19; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.0₃0906"; E = "  +6.11%  " },
    @{ Row = 20; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "72.28"; E = "  +0.04%  " },
    @{ Row = 21; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "6.07"; E = "  -1.25%  " },
    @{ Row = 22; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "232.58"; E = "  +1.58%  " },
    @{ Row = 23; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.00"; E = "  -0.01%  " },
    @{ Row = 24; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "2.39"; E = "  +3.51%  " },
    @{ Row = 25; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "2.35"; E = "  -0.97%  " },
    @{ Row = 26; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "9.68"; E = "  +0.18%  " },
    @{ Row = 27; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "171.79"; E = "  -0.13%  " },
    @{ Row = 28; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "0.140"; E = "  +1.68%  " },
    @{ Row = 29; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "1.45"; E = "  +2.19%  " },
    @{ Row = 30; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "20.06"; E = "  +1.70%  " },
    @{ Row = 31; B = "WEMIXToken"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D = "2.74"; E = "  +4.66%  " },
    @{ Row = 32; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.123"; E = "  +0.96%  " },
    @{ Row = 33; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "4.57"; E = "  -1.82%  " },
    @{ Row = 34; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "4.72"; E = "  -1.94%  " },
    @{ Row = 35; B = "THORChain"; C = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"; D = "7.04"; E = "  -0.70%  " },
    @{ Row = 36; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.0623"; E = "  -0.06%  " },
    @{ Row = 37; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "3.91"; E = "  +8.10%  " },
    @{ Row = 38; B = "LidoDAOToken"; C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D = "2.43"; E = "  -0.13%  " },
    @{ Row = 39; B = "FTXToken"; C = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D = "5.04"; E = "  +19.86%  " },
    @{ Row = 40; B = "BinanceUSD"; C = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D = "1.00"; E = "  +0.13%  " },
    @{ Row = 41; B = "Aave"; C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D = "103.41"; E = "  -1.18%  " },
    @{ Row = 42; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.0229"; E = "  -0.98%  " },
    @{ Row = 43; B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D = "8.29"; E = "  +4.85%  " },
    @{ Row = 44; B = "Maker"; C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D = "1.518.39"; E = "  -1.30%  " },
    @{ Row = 45; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "17.40"; E = "  -3.19%  " },
    @{ Row = 46; B = "TrustWalletToken"; C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D = "1.22"; E = "  +2.10%  " },
    @{ Row = 47; B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.0928"; E = "  -0.61%  " },
    @{ Row = 48; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "1.10"; E = "  -0.10%  " },
    @{ Row = 49; B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "2.80"; E = "  -0.48%  " },
    @{ Row = 50; B = "TerraClassic"; C = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"; D = "0.000195"; E = "  +31.91%  " },
    @{ Row = 51; B = "RocketPoolETH"; C = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; D = "2.404.31"; E = "  +1.55%  " }
)

foreach ($item in $data) {
    $r = $item.Row

    $ws.Cells.Item($r, 2).Value = $item.B   # Coin
    $ws.Cells.Item($r, 3).Value = $item.C   # Link
    $ws.Cells.Item($r, 5).Value = $item.E   # Volume(1h)

    # Column D (Price) holds text that often looks numeric (e.g. "228.00",
    # "40.019.13", "0.630"); a plain .Value assignment would let Excel's
    # auto-detection coerce it to a number and silently drop formatting
    # such as trailing zeros. Force text entry, then restore the cell's
    # original (unstyled) format so no stray style gets introduced.
    $priceCell = $ws.Cells.Item($r, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $item.D
    $priceCell.Style = $ws.Cells.Item($r, 2).Style
}
